$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the "arquivo" value for Tonho's row to the new file name
$ws.Range("C2").Value = "transferir.png"

# Move the active selection from A5 to C5
$ws.Range("C5").Select()
